# Update cryptos list values (price/volume refresh) per the Nov 5 2024 GitHub Actions run.
# Also swaps the SuiNetwork/Dai rows (22 and 23), matching upstream coinranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.056.14'
$ws.Cells.Item(2, 5).Value = '  -1.59%  '
$ws.Cells.Item(3, 4).Value = '2.413.59'
$ws.Cells.Item(3, 5).Value = '  -2.56%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '554.93'
$ws.Cells.Item(5, 5).Value = '  -1.52%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '158.92'
$ws.Cells.Item(6, 5).Value = '  -2.74%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$ws.Cells.Item(8, 5).Value = '  -0.44%  '
$ws.Cells.Item(9, 5).Value = '  +6.42%  '
$ws.Cells.Item(10, 5).Value = '  -1.75%  '
$ws.Cells.Item(11, 5).Value = '  -1.31%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '4.64'
$ws.Cells.Item(12, 5).Value = '  -4.83%  '
$ws.Cells.Item(13, 4).Value = '67.942.69'
$ws.Cells.Item(13, 5).Value = '  -1.58%  '
$ws.Cells.Item(14, 4).Value = '2.852.72'
$ws.Cells.Item(14, 5).Value = '  -1.83%  '
$ws.Cells.Item(15, 5).Value = '  +1.53%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '22.80'
$ws.Cells.Item(16, 5).Value = '  -4.09%  '
$ws.Cells.Item(17, 4).Value = '2.410.45'
$ws.Cells.Item(17, 5).Value = '  -2.72%  '
$ws.Cells.Item(18, 5).Value = '  -4.17%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '331.01'
$ws.Cells.Item(19, 5).Value = '  -2.80%  '
$ws.Cells.Item(20, 5).Value = '  -3.10%  '
$ws.Cells.Item(21, 5).Value = '  -0.53%  '
$ws.Cells.Item(22, 2).Value = 'SuiNetwork'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '1.88'
$ws.Cells.Item(22, 5).Value = '  -1.89%  '
$ws.Cells.Item(23, 2).Value = 'Dai'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.00'
$ws.Cells.Item(23, 5).Value = '  -0.02%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '66.15'
$ws.Cells.Item(24, 5).Value = '  -1.89%  '
$ws.Cells.Item(25, 5).Value = '  -1.80%  '
$ws.Cells.Item(26, 4).Value = '2.535.45'
$ws.Cells.Item(26, 5).Value = '  -2.88%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '8.16'
$ws.Cells.Item(27, 5).Value = '  -1.35%  '
$ws.Cells.Item(28, 5).Value = '  -2.44%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.08'
$ws.Cells.Item(29, 5).Value = '  -2.06%  '
$ws.Cells.Item(30, 5).Value = '  +0.11%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '420.63'
$ws.Cells.Item(31, 5).Value = '  -3.48%  '
$ws.Cells.Item(32, 5).Value = '  -1.88%  '
$ws.Cells.Item(33, 5).Value = '  -2.17%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '158.37'
$ws.Cells.Item(34, 5).Value = '  +0.54%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '19.02'
$ws.Cells.Item(35, 5).Value = '  -0.23%  '
$ws.Cells.Item(36, 5).Value = '  -0.01%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '17.80'
$ws.Cells.Item(37, 5).Value = '  -0.47%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.105'
$ws.Cells.Item(39, 5).Value = '  -2.48%  '
$ws.Cells.Item(40, 5).Value = '  -4.05%  '
$ws.Cells.Item(41, 5).Value = '  -1.08%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.07'
$ws.Cells.Item(42, 5).Value = '  -1.65%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '132.34'
$ws.Cells.Item(43, 5).Value = '  -1.17%  '
$ws.Cells.Item(44, 5).Value = '  -1.60%  '
$ws.Cells.Item(45, 5).Value = '  -5.77%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0712'
$ws.Cells.Item(46, 5).Value = '  -0.82%  '
$ws.Cells.Item(47, 5).Value = '  -1.96%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.554'
$ws.Cells.Item(48, 5).Value = '  -1.93%  '
$ws.Cells.Item(49, 5).Value = '  -0.16%  '
$ws.Cells.Item(51, 5).Value = '  -3.54%  '
